# January-2021-Collection.xlsx — "data updated till 12Jan 8AM"
#
# Adds the Jan-11-2021 collection column (Q) figures for a batch of
# retailers, and records a new collection agent "Abhimanyu" as the alias
# for the retailer in row 4 (which previously had no alias recorded).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New alias name for the retailer on row 4 ("Abhimanyu")
$ws.Range("B4").Value = "Abhimanyu"

# New collection entries for 11-Jan-2021 (column Q)
$ws.Range("Q4").Value  = 1500
$ws.Range("Q5").Value  = 1000
$ws.Range("Q6").Value  = 3000
$ws.Range("Q14").Value = 2500
$ws.Range("Q17").Value = 2500
$ws.Range("Q20").Value = 1000
$ws.Range("Q22").Value = 1000
$ws.Range("Q25").Value = 3000
$ws.Range("Q31").Value = 1000
$ws.Range("Q32").Value = 2000
$ws.Range("Q35").Value = 3000
$ws.Range("Q41").Value = 6000
$ws.Range("Q44").Value = 3000
$ws.Range("Q46").Value = 3000
$ws.Range("Q48").Value = 3000
$ws.Range("Q51").Value = 1000
$ws.Range("Q52").Value = 1000
$ws.Range("Q53").Value = 1000
$ws.Range("Q62").Value = 2000
$ws.Range("Q65").Value = 10000
$ws.Range("Q69").Value = 500
$ws.Range("Q71").Value = 2000

# Leave the view scrolled/selected where the author last left off.
$excel.Goto($ws.Range("Q98"))
